# Auto-generated Word COM-interop edit script.
# Applies: heading font-size bump (sz/szCs 28) on the four 'TitoloN' style
# headings, and spell-check proofErr run-splitting around GiocoPadel/
# padeleur/Padeleur/pedaleur occurrences, matching the target diff.

$d = $word.ActiveDocument

function Set-HeadingSize($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Font.Size = 14
    $r.Font.SizeBi = 14
}

function Set-ParagraphXml($paraIndex, [string]$b64) {
    $bytes = [System.Convert]::FromBase64String($b64)
    $xml = [System.Text.Encoding]::UTF8.GetString($bytes)
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range.Duplicate
    $r.InsertXML($xml)
}

# --- 1) Heading font sizes: '1 Introduzione', '2 Descrizione',
#        '3 Descrizione delle parti interessate',
#        '4 Riepilogo delle caratteristiche del sistema' ---
Set-HeadingSize 2
Set-HeadingSize 4
Set-HeadingSize 8
Set-HeadingSize 11

# --- 2) Paragraphs whose runs get split w/ w:proofErr spellcheck markers ---
Set-ParagraphXml 3 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpwU3R5bGUgdzp2YWw9Ik5lc3N1bmFzcGF6aWF0dXJhIi8+PHc6amMgdzp2YWw9ImJvdGgiLz48dzpyUHI+PHc6c3ogdzp2YWw9IjI0Ii8+PHc6c3pDcyB3OnZhbD0iMjQiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6c3ogdzp2YWw9IjI0Ii8+PHc6c3pDcyB3OnZhbD0iMjQiLz48L3c6clByPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+UHJldmVkaWFtbyBsYSByZWFsaXp6YXppb25lIGRpIHVu4oCZYXBwbGljYXppb25lIGRpIGdlc3Rpb25lIGRpIHVu4oCZYXNzb2NpYXppb25lIHNwb3J0aXZhIGRpIHBhZGVsLCBjaGlhbWF0YSA8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9InNwZWxsU3RhcnQiLz48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dD5HaW9jb1BhZGVsPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbEVuZCIvPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPi4gTOKAmW9iaWV0dGl2byBkZWxs4oCZYXBwbGljYXppb25lIMOoIHF1ZWxsbyBkaSBnZXN0aXJlIGkgdHJlIGNhbXBpIGRhIHBhZGVsIGFwcGFydGVuZW50aSBhbGzigJlhc3NvY2lhemlvbmUgY29uIGxhIHJlbGF0aXZhIHByZW5vdGF6aW9uZSBkZWkgY2FtcGkgZGEgcGFydGUgZGVpIDwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxTdGFydCIvPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0PnBhZGVsZXVyPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbEVuZCIvPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0PiwgZGVsbOKAmWV2ZW50dWFsZSBhdHRyZXp6YXR1cmEgbm9uY2jDqSBkZWxsYSBnZXN0aW9uZSBkaSBwb2xpdGljYSBkZWwgcHJlenpvLjwvdzp0PjwvdzpyPjwvdzpwPg=="
Set-ParagraphXml 5 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpqYyB3OnZhbD0iYm90aCIvPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PC93OnBQcj48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5JbCBzb2Z0d2FyZSA8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9InNwZWxsU3RhcnQiLz48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dD5HaW9jb1BhZGVsPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbEVuZCIvPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBzaSBwb25lIGzigJlvYmlldHRpdm8gZGkgc29zdGl0dWlyZSBpIHRyYWRpemlvbmFsaSBtZXRvZGkgZGkgZ2VzdGlvbmUgbm9uIGF1dG9tYXRpenphdGEgZGVpIGNhbXBpIHBlciBmb3JuaXJlIGFuY2hlIGRhdGkgc3RhdGlzdGljaSBhbCBmaW5lIGRpIG1hc3NpbWl6emFyZSBpIGd1YWRhZ25pIGFuZGFuZG8gYSBtb2RpZmljYXJlIGkgcHJlenppIGRlaSBjYW1waSBkaSBwYWRlbC4gPC93OnQ+PC93OnI+PHc6cj48dzpyUHI+PHc6c3ogdzp2YWw9IjI0Ii8+PHc6c3pDcyB3OnZhbD0iMjQiLz48L3c6clByPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IDwvdzp0PjwvdzpyPjwvdzpwPg=="
Set-ParagraphXml 10 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpwU3R5bGUgdzp2YWw9IlBhcmFncmFmb2VsZW5jbyIvPjx3Om51bVByPjx3Omlsdmwgdzp2YWw9IjAiLz48dzpudW1JZCB3OnZhbD0iNCIvPjwvdzpudW1Qcj48dzpqYyB3OnZhbD0iYm90aCIvPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PC93OnBQcj48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dD5VdGVudGkgKDwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxTdGFydCIvPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0PlBhZGVsZXVyPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbEVuZCIvPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0Pik6IHBlcnNvbmUgY2hlIG5lY2Vzc2l0YW5vIHVuIHBvcnRhbGUgcGVyIGxhIHByZW5vdGF6aW9uZSBkZWkgY2FtcGk8L3c6dD48L3c6cj48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gZWQgZXZlbnR1YWxlIGF0dHJlenphdHVyYTwvdzp0PjwvdzpyPjx3OnI+PHc6clByPjx3OnN6IHc6dmFsPSIyNCIvPjx3OnN6Q3Mgdzp2YWw9IjI0Ii8+PC93OnJQcj48dzp0Pi48L3c6dD48L3c6cj48L3c6cD4="
Set-ParagraphXml 13 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpwU3R5bGUgdzp2YWw9IlBhcmFncmFmb2VsZW5jbyIvPjx3Om51bVByPjx3Omlsdmwgdzp2YWw9IjAiLz48dzpudW1JZCB3OnZhbD0iNiIvPjwvdzpudW1Qcj48dzpqYyB3OnZhbD0iYm90aCIvPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PC93OnBQcj48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5SZWdpc3RyYXppb25lIGRlaSA8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9InNwZWxsU3RhcnQiLz48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dD5wZWRhbGV1cjwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxFbmQiLz48dzpyPjx3OnJQcj48dzpzeiB3OnZhbD0iMjQiLz48dzpzekNzIHc6dmFsPSIyNCIvPjwvdzpyUHI+PHc6dD47PC93OnQ+PC93OnI+PC93OnA+"

Write-Output "edit applied"
